# Adicion del IMT con una prueba para el clasificador
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 45870
$ws.Range("B2").Value = -0.3
$ws.Range("C2").Value = 'ONI'
$ws.Range("D2").Value = 'Índice Niño Oceánico: Media móvil de 3 meses de las anomalías de la TSM ERSST.v5 en la región Niño 3.4 (5°N-5°S, 120°-170°W) Calculada a partir del ERSST V5 (en NOAA/CPC).'
$ws.Range("E2").Value = '°C'
$ws.Range("F2").Value = 'Neutra'
$ws.Range("G2").Value = 'Esta fase se caracteriza porque las anomalías de TSM son inferiores a 0.5 °C y superiores a -0.5 °C'
$ws.Range("H2").Value = 'Neutro'
$ws.Range("I2").Value = 'Condiciones neutras'
$ws.Range("J2").Value = 'Neutro'

# Row 3
$ws.Range("A3").Value = 45901
$ws.Range("B3").Value = -0.4
$ws.Range("C3").Value = 'ONI'
$ws.Range("D3").Value = 'Índice Niño Oceánico: Media móvil de 3 meses de las anomalías de la TSM ERSST.v5 en la región Niño 3.4 (5°N-5°S, 120°-170°W) Calculada a partir del ERSST V5 (en NOAA/CPC).'
$ws.Range("E3").Value = '°C'
$ws.Range("F3").Value = 'Neutra'
$ws.Range("G3").Value = 'Esta fase se caracteriza porque las anomalías de TSM son inferiores a 0.5 °C y superiores a -0.5 °C'
$ws.Range("H3").Value = 'Neutro'
$ws.Range("I3").Value = 'Condiciones neutras'
$ws.Range("J3").Value = 'Neutro'

# Row 4
$ws.Range("A4").Value = 45901
$ws.Range("B4").Value = -0.4
$ws.Range("C4").Value = 'Niño 1+2'
$ws.Range("D4").Value = 'Índice Niño 1+2: representa las anomalías mensuales de la temperatura superficial del mar (TSM) en la región más oriental del Pacífico ecuatorial, delimitada entre los 0°–10°S y 80°W–90°W, frente a las costas de Perú y Ecuador. Calculada a partir del ERSST V5 (en NOAA/CPC).'
$ws.Range("E4").Value = '°C'
$ws.Range("F4").Value = 'Neutra'
$ws.Range("G4").Value = 'Esta fase se caracteriza porque las anomalías de TSM son inferiores a 0.5 °C y superiores a -0.5 °C'
$ws.Range("H4").Value = 'Neutro'
$ws.Range("I4").Value = 'Condiciones neutras'
$ws.Range("J4").Value = 'Neutro'

# Row 5
$ws.Range("A5").Value = 45931
$ws.Range("B5").Value = -0.3
$ws.Range("C5").Value = 'Niño 1+2'
$ws.Range("D5").Value = 'Índice Niño 1+2: representa las anomalías mensuales de la temperatura superficial del mar (TSM) en la región más oriental del Pacífico ecuatorial, delimitada entre los 0°–10°S y 80°W–90°W, frente a las costas de Perú y Ecuador. Calculada a partir del ERSST V5 (en NOAA/CPC).'
$ws.Range("E5").Value = '°C'
$ws.Range("F5").Value = 'Neutra'
$ws.Range("G5").Value = 'Esta fase se caracteriza porque las anomalías de TSM son inferiores a 0.5 °C y superiores a -0.5 °C'
$ws.Range("H5").Value = 'Neutro'
$ws.Range("I5").Value = 'Condiciones neutras'
$ws.Range("J5").Value = 'Neutro'

# Row 6
$ws.Range("A6").Value = 45901
$ws.Range("B6").Value = -0.5
$ws.Range("C6").Value = 'Niño 3'
$ws.Range("D6").Value = 'Índice Niño 3: El índice Niño 3 corresponde a las anomalías mensuales de la temperatura superficial del mar (TSM) en la región del Pacífico ecuatorial comprendida entre los 5°N–5°S y 90°W–150°W. Calculada a partir del ERSST V5 (en NOAA/CPC).'
$ws.Range("E6").Value = '°C'
$ws.Range("F6").Value = 'Fría'
$ws.Range("G6").Value = 'Esta fase se caracteriza porque las anomalías de TSM en la región 3 son inferiores a -0.5 °C'
$ws.Range("H6").Value = 'Neutro'
$ws.Range("I6").Value = 'Condiciones neutras'
$ws.Range("J6").Value = 'Neutro'

# Row 7
$ws.Range("A7").Value = 45931
$ws.Range("B7").Value = -0.4
$ws.Range("C7").Value = 'Niño 3'
$ws.Range("D7").Value = 'Índice Niño 3: El índice Niño 3 corresponde a las anomalías mensuales de la temperatura superficial del mar (TSM) en la región del Pacífico ecuatorial comprendida entre los 5°N–5°S y 90°W–150°W. Calculada a partir del ERSST V5 (en NOAA/CPC).'
$ws.Range("E7").Value = '°C'
$ws.Range("F7").Value = 'Neutra'
$ws.Range("G7").Value = 'Esta fase se caracteriza porque las anomalías de TSM son inferiores a 0.5 °C y superiores a -0.5 °C'
$ws.Range("H7").Value = 'Neutro'
$ws.Range("I7").Value = 'Condiciones neutras'
$ws.Range("J7").Value = 'Neutro'

# Row 8
$ws.Range("A8").Value = 45901
$ws.Range("B8").Value = -0.5
$ws.Range("C8").Value = 'Niño 3.4'
$ws.Range("D8").Value = 'Índice Niño 3.4: El índice Niño 3.4 mide las anomalías mensuales de la temperatura superficial del mar (TSM) en la región comprendida entre los 5°N–5°S y 120°W–170°W del Pacífico central ecuatorial. Calculada a partir del ERSST V5 (en NOAA/CPC).'
$ws.Range("E8").Value = '°C'
$ws.Range("F8").Value = 'Fría'
$ws.Range("G8").Value = 'Esta fase se caracteriza porque las anomalías de TSM en la región 3.4 son inferiores a -0.5 °C'
$ws.Range("H8").Value = 'Neutro'
$ws.Range("I8").Value = 'Condiciones neutras'
$ws.Range("J8").Value = 'Neutro'

# Row 9
$ws.Range("A9").Value = 45931
$ws.Range("B9").Value = -0.5
$ws.Range("C9").Value = 'Niño 3.4'
$ws.Range("D9").Value = 'Índice Niño 3.4: El índice Niño 3.4 mide las anomalías mensuales de la temperatura superficial del mar (TSM) en la región comprendida entre los 5°N–5°S y 120°W–170°W del Pacífico central ecuatorial. Calculada a partir del ERSST V5 (en NOAA/CPC).'
$ws.Range("E9").Value = '°C'
$ws.Range("F9").Value = 'Fría'
$ws.Range("G9").Value = 'Esta fase se caracteriza porque las anomalías de TSM en la región 3.4 son inferiores a -0.5 °C'
$ws.Range("H9").Value = 'Neutro'
$ws.Range("I9").Value = 'Condiciones neutras'
$ws.Range("J9").Value = 'Neutro'

# Row 10
$ws.Range("A10").Value = 45901
$ws.Range("B10").Value = -0.3
$ws.Range("C10").Value = 'Niño 4'
$ws.Range("D10").Value = 'Índice Niño 4: El índice Niño 4 representa las anomalías mensuales de la temperatura superficial del mar (TSM) en la región del Pacífico ecuatorial occidental, delimitada entre los 5°N–5°S y 160°E–150°W. Calculada a partir del ERSST V5 (en NOAA/CPC).'
$ws.Range("E10").Value = '°C'
$ws.Range("F10").Value = 'Neutra'
$ws.Range("G10").Value = 'Esta fase se caracteriza porque las anomalías de TSM son inferiores a 0.5 °C y superiores a -0.5 °C'
$ws.Range("H10").Value = 'Neutro'
$ws.Range("I10").Value = 'Condiciones neutras'
$ws.Range("J10").Value = 'Neutro'

# Row 11
$ws.Range("A11").Value = 45931
$ws.Range("B11").Value = -0.4
$ws.Range("C11").Value = 'Niño 4'
$ws.Range("D11").Value = 'Índice Niño 4: El índice Niño 4 representa las anomalías mensuales de la temperatura superficial del mar (TSM) en la región del Pacífico ecuatorial occidental, delimitada entre los 5°N–5°S y 160°E–150°W. Calculada a partir del ERSST V5 (en NOAA/CPC).'
$ws.Range("E11").Value = '°C'
$ws.Range("F11").Value = 'Neutra'
$ws.Range("G11").Value = 'Esta fase se caracteriza porque las anomalías de TSM son inferiores a 0.5 °C y superiores a -0.5 °C'
$ws.Range("H11").Value = 'Neutro'
$ws.Range("I11").Value = 'Condiciones neutras'
$ws.Range("J11").Value = 'Neutro'

# Row 12
$ws.Range("A12").Value = 45901
$ws.Range("B12").Value = 0.1
$ws.Range("C12").Value = 'SOI'
$ws.Range("D12").Value = 'Southern Oscillation Index: El Índice de la Oscilación del Sur es un indicador climático que mide la diferencia de presión atmosférica a nivel del mar entre dos estaciones del Pacífico tropical: Tahití (Polinesia Francesa) y Darwin (Australia). Calculada a partir del ERSST V5 (en NOAA/CPC https://www.psl.noaa.gov/data/timeseries/month/DS/SOI/).'
$ws.Range("E12").Value = 'dmLess'
$ws.Range("F12").Value = 'Fría'
$ws.Range("G12").Value = 'Esta fase se caracteriza por presiones más altas en Tahití y más bajas en Darwin, típicas de La Niña (SOI positivo)'
$ws.Range("H12").Value = 'Neutro'
$ws.Range("I12").Value = 'Este evento se caracteriza porque el valor del índice para el mes es cero'
$ws.Range("J12").Value = 'No aplicable'

# Row 13
$ws.Range("A13").Value = 45931
$ws.Range("B13").Value = 1.9
$ws.Range("C13").Value = 'SOI'
$ws.Range("D13").Value = 'Southern Oscillation Index: El Índice de la Oscilación del Sur es un indicador climático que mide la diferencia de presión atmosférica a nivel del mar entre dos estaciones del Pacífico tropical: Tahití (Polinesia Francesa) y Darwin (Australia). Calculada a partir del ERSST V5 (en NOAA/CPC https://www.psl.noaa.gov/data/timeseries/month/DS/SOI/).'
$ws.Range("E13").Value = 'dmLess'
$ws.Range("F13").Value = 'Fría'
$ws.Range("G13").Value = 'Esta fase se caracteriza por presiones más altas en Tahití y más bajas en Darwin, típicas de La Niña (SOI positivo)'
$ws.Range("H13").Value = 'Neutro'
$ws.Range("I13").Value = 'Este evento se caracteriza porque el valor del índice para el mes es cero'
$ws.Range("J13").Value = 'No aplicable'

# Clear trailing rows that no longer have data
$ws.Range("A14:J17").ClearContents()

# Column D (index_description) is no longer hidden
$ws.Columns(4).Hidden = $false

# Restore the selection left by the editor
$ws.Range("D17").Select()
